# Rewrite the "Test_format_trace" sheet with the refreshed rail-car trace
# report (17 events instead of 9; color-code-driven relabeling of the
# Initial/Location/Event columns per the updated trace pull).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wipe the old A1:O11 report body (and anything already in A1:O19) so
# stale cells do not linger once the new table is narrower/shorter in a column.
$ws.Range("A1:O19").ClearContents()

$ws.Cells.Item(1,1).Value = "Description unknown, completed 10/24/2023 07:54:28 EDT, by WPJTOWN1.The search returned: 17 events."

$ws.Cells.Item(2,1).Value = "Initial"
$ws.Cells.Item(2,2).Value = "Number"
$ws.Cells.Item(2,3).Value = "Location City"
$ws.Cells.Item(2,4).Value = "State"
$ws.Cells.Item(2,5).Value = "Month"
$ws.Cells.Item(2,6).Value = "Day"
$ws.Cells.Item(2,7).Value = "Time"
$ws.Cells.Item(2,8).Value = "Event"
$ws.Cells.Item(2,9).Value = "Train ID"
$ws.Cells.Item(2,10).Value = "Destination City"
$ws.Cells.Item(2,11).Value = "State"
$ws.Cells.Item(2,12).Value = "Gross Weight"
$ws.Cells.Item(2,13).Value = "Tare Weight"
$ws.Cells.Item(2,14).Value = "Net Weight"
$ws.Cells.Item(2,15).Value = "Car_no"

$ws.Cells.Item(3,1).Value = "TCIX"
$ws.Cells.Item(3,2).Value = 250760
$ws.Cells.Item(3,3).Value = "CHEYENNE"
$ws.Cells.Item(3,4).Value = "WY"
$ws.Cells.Item(3,5).Value = 10
$ws.Cells.Item(3,6).Value = 23
$ws.Cells.Item(3,7).Value = 2346
$ws.Cells.Item(3,8).Value = "Arrive In-Transit"
$ws.Cells.Item(3,9).Value = "MNPDV2"
$ws.Cells.Item(3,10).Value = "JOHNSTOWN"
$ws.Cells.Item(3,11).Value = "CO"
$ws.Cells.Item(3,12).Value = 170000
$ws.Cells.Item(3,13).Value = 0
$ws.Cells.Item(3,14).Value = 170000
$ws.Cells.Item(3,15).Value = "TCIX250760"

$ws.Cells.Item(4,1).Value = "TCIX"
$ws.Cells.Item(4,2).Value = 250766
$ws.Cells.Item(4,3).Value = "CHEYENNE"
$ws.Cells.Item(4,4).Value = "WY"
$ws.Cells.Item(4,5).Value = 10
$ws.Cells.Item(4,6).Value = 23
$ws.Cells.Item(4,7).Value = 2346
$ws.Cells.Item(4,8).Value = "Arrive In-Transit"
$ws.Cells.Item(4,9).Value = "MNPDV2"
$ws.Cells.Item(4,10).Value = "JOHNSTOWN"
$ws.Cells.Item(4,11).Value = "CO"
$ws.Cells.Item(4,12).Value = 170000
$ws.Cells.Item(4,13).Value = 0
$ws.Cells.Item(4,14).Value = 170000
$ws.Cells.Item(4,15).Value = "TCIX250766"

$ws.Cells.Item(5,1).Value = "TILX"
$ws.Cells.Item(5,2).Value = 252319
$ws.Cells.Item(5,3).Value = "CRESTON"
$ws.Cells.Item(5,4).Value = "WY"
$ws.Cells.Item(5,5).Value = 10
$ws.Cells.Item(5,6).Value = 24
$ws.Cells.Item(5,7).Value = 241
$ws.Cells.Item(5,8).Value = "Departure"
$ws.Cells.Item(5,9).Value = "MWCNP2"
$ws.Cells.Item(5,10).Value = "JOHNSTOWN"
$ws.Cells.Item(5,11).Value = "CO"
$ws.Cells.Item(5,12).Value = 180150
$ws.Cells.Item(5,13).Value = 0
$ws.Cells.Item(5,14).Value = 180150
$ws.Cells.Item(5,15).Value = "TILX252319"

$ws.Cells.Item(6,1).Value = "UTLX"
$ws.Cells.Item(6,2).Value = 645561
$ws.Cells.Item(6,3).Value = "CRESTON"
$ws.Cells.Item(6,4).Value = "WY"
$ws.Cells.Item(6,5).Value = 10
$ws.Cells.Item(6,6).Value = 24
$ws.Cells.Item(6,7).Value = 241
$ws.Cells.Item(6,8).Value = "Departure"
$ws.Cells.Item(6,9).Value = "MWCNP2"
$ws.Cells.Item(6,10).Value = "JOHNSTOWN"
$ws.Cells.Item(6,11).Value = "CO"
$ws.Cells.Item(6,12).Value = 180000
$ws.Cells.Item(6,13).Value = 0
$ws.Cells.Item(6,14).Value = 180000
$ws.Cells.Item(6,15).Value = "UTLX645561"

$ws.Cells.Item(7,1).Value = "TILX"
$ws.Cells.Item(7,2).Value = 252283
$ws.Cells.Item(7,3).Value = "GREELEY"
$ws.Cells.Item(7,4).Value = "CO"
$ws.Cells.Item(7,5).Value = 10
$ws.Cells.Item(7,6).Value = 22
$ws.Cells.Item(7,7).Value = 1506
$ws.Cells.Item(7,8).Value = "Junction Delivery"
$ws.Cells.Item(7,9).Value = "GWR"
$ws.Cells.Item(7,10).Value = "JOHNSTOWN"
$ws.Cells.Item(7,11).Value = "CO"
$ws.Cells.Item(7,12).Value = 178200
$ws.Cells.Item(7,13).Value = 0
$ws.Cells.Item(7,14).Value = 178200
$ws.Cells.Item(7,15).Value = "TILX252283"

$ws.Cells.Item(8,1).Value = "TCIX"
$ws.Cells.Item(8,2).Value = 256500
$ws.Cells.Item(8,3).Value = "JOHNSTOWN"
$ws.Cells.Item(8,4).Value = "CO"
$ws.Cells.Item(8,5).Value = 10
$ws.Cells.Item(8,6).Value = 9
$ws.Cells.Item(8,7).Value = 1330
$ws.Cells.Item(8,8).Value = "Placed Actual"
$ws.Cells.Item(8,10).Value = "JOHNSTOWN"
$ws.Cells.Item(8,11).Value = "CO"
$ws.Cells.Item(8,12).Value = 179700
$ws.Cells.Item(8,13).Value = 0
$ws.Cells.Item(8,14).Value = 179700
$ws.Cells.Item(8,15).Value = "TCIX256500"

$ws.Cells.Item(9,1).Value = "TCIX"
$ws.Cells.Item(9,2).Value = 256424
$ws.Cells.Item(9,3).Value = "JOHNSTOWN"
$ws.Cells.Item(9,4).Value = "CO"
$ws.Cells.Item(9,5).Value = 10
$ws.Cells.Item(9,6).Value = 12
$ws.Cells.Item(9,7).Value = 1626
$ws.Cells.Item(9,8).Value = "Placed Actual"
$ws.Cells.Item(9,10).Value = "JOHNSTOWN"
$ws.Cells.Item(9,11).Value = "CO"
$ws.Cells.Item(9,12).Value = 179650
$ws.Cells.Item(9,13).Value = 0
$ws.Cells.Item(9,14).Value = 179650
$ws.Cells.Item(9,15).Value = "TCIX256424"

$ws.Cells.Item(10,1).Value = "TCIX"
$ws.Cells.Item(10,2).Value = 250780
$ws.Cells.Item(10,3).Value = "JOHNSTOWN"
$ws.Cells.Item(10,4).Value = "CO"
$ws.Cells.Item(10,5).Value = 10
$ws.Cells.Item(10,6).Value = 12
$ws.Cells.Item(10,7).Value = 1626
$ws.Cells.Item(10,8).Value = "Placed Actual"
$ws.Cells.Item(10,10).Value = "JOHNSTOWN"
$ws.Cells.Item(10,11).Value = "CO"
$ws.Cells.Item(10,12).Value = 179600
$ws.Cells.Item(10,13).Value = 0
$ws.Cells.Item(10,14).Value = 179600
$ws.Cells.Item(10,15).Value = "TCIX250780"

$ws.Cells.Item(11,1).Value = "UTLX"
$ws.Cells.Item(11,2).Value = 645560
$ws.Cells.Item(11,3).Value = "JOHNSTOWN"
$ws.Cells.Item(11,4).Value = "CO"
$ws.Cells.Item(11,5).Value = 10
$ws.Cells.Item(11,6).Value = 17
$ws.Cells.Item(11,7).Value = 1415
$ws.Cells.Item(11,8).Value = "Placed Actual"
$ws.Cells.Item(11,10).Value = "JOHNSTOWN"
$ws.Cells.Item(11,11).Value = "CO"
$ws.Cells.Item(11,12).Value = 179750
$ws.Cells.Item(11,13).Value = 0
$ws.Cells.Item(11,14).Value = 179750
$ws.Cells.Item(11,15).Value = "UTLX645560"

$ws.Cells.Item(12,1).Value = "UTLX"
$ws.Cells.Item(12,2).Value = 645570
$ws.Cells.Item(12,3).Value = "JOHNSTOWN"
$ws.Cells.Item(12,4).Value = "CO"
$ws.Cells.Item(12,5).Value = 10
$ws.Cells.Item(12,6).Value = 17
$ws.Cells.Item(12,7).Value = 1415
$ws.Cells.Item(12,8).Value = "Placed Actual"
$ws.Cells.Item(12,10).Value = "JOHNSTOWN"
$ws.Cells.Item(12,11).Value = "CO"
$ws.Cells.Item(12,12).Value = 179550
$ws.Cells.Item(12,13).Value = 0
$ws.Cells.Item(12,14).Value = 179550
$ws.Cells.Item(12,15).Value = "UTLX645570"

$ws.Cells.Item(13,1).Value = "TCIX"
$ws.Cells.Item(13,2).Value = 258654
$ws.Cells.Item(13,3).Value = "JOHNSTOWN"
$ws.Cells.Item(13,4).Value = "CO"
$ws.Cells.Item(13,5).Value = 10
$ws.Cells.Item(13,6).Value = 17
$ws.Cells.Item(13,7).Value = 1415
$ws.Cells.Item(13,8).Value = "Placed Actual"
$ws.Cells.Item(13,10).Value = "JOHNSTOWN"
$ws.Cells.Item(13,11).Value = "CO"
$ws.Cells.Item(13,12).Value = 173850
$ws.Cells.Item(13,13).Value = 0
$ws.Cells.Item(13,14).Value = 173850
$ws.Cells.Item(13,15).Value = "TCIX258654"

$ws.Cells.Item(14,1).Value = "TCIX"
$ws.Cells.Item(14,2).Value = 256434
$ws.Cells.Item(14,3).Value = "JOHNSTOWN"
$ws.Cells.Item(14,4).Value = "CO"
$ws.Cells.Item(14,5).Value = 10
$ws.Cells.Item(14,6).Value = 17
$ws.Cells.Item(14,7).Value = 1415
$ws.Cells.Item(14,8).Value = "Placed Actual"
$ws.Cells.Item(14,10).Value = "JOHNSTOWN"
$ws.Cells.Item(14,11).Value = "CO"
$ws.Cells.Item(14,12).Value = 179950
$ws.Cells.Item(14,13).Value = 0
$ws.Cells.Item(14,14).Value = 179950
$ws.Cells.Item(14,15).Value = "TCIX256434"

$ws.Cells.Item(15,1).Value = "TCIX"
$ws.Cells.Item(15,2).Value = 256419
$ws.Cells.Item(15,3).Value = "JOHNSTOWN"
$ws.Cells.Item(15,4).Value = "CO"
$ws.Cells.Item(15,5).Value = 10
$ws.Cells.Item(15,6).Value = 18
$ws.Cells.Item(15,7).Value = 1602
$ws.Cells.Item(15,8).Value = "Placed Actual"
$ws.Cells.Item(15,10).Value = "JOHNSTOWN"
$ws.Cells.Item(15,11).Value = "CO"
$ws.Cells.Item(15,12).Value = 180100
$ws.Cells.Item(15,13).Value = 0
$ws.Cells.Item(15,14).Value = 180100
$ws.Cells.Item(15,15).Value = "TCIX256419"

$ws.Cells.Item(16,1).Value = "UTLX"
$ws.Cells.Item(16,2).Value = 669035
$ws.Cells.Item(16,3).Value = "JOHNSTOWN"
$ws.Cells.Item(16,4).Value = "CO"
$ws.Cells.Item(16,5).Value = 10
$ws.Cells.Item(16,6).Value = 23
$ws.Cells.Item(16,7).Value = 1457
$ws.Cells.Item(16,8).Value = "Placed Actual"
$ws.Cells.Item(16,10).Value = "JOHNSTOWN"
$ws.Cells.Item(16,11).Value = "CO"
$ws.Cells.Item(16,12).Value = 179850
$ws.Cells.Item(16,13).Value = 0
$ws.Cells.Item(16,14).Value = 179850
$ws.Cells.Item(16,15).Value = "UTLX669035"

$ws.Cells.Item(17,1).Value = "UTLX"
$ws.Cells.Item(17,2).Value = 669036
$ws.Cells.Item(17,3).Value = "JOHNSTOWN"
$ws.Cells.Item(17,4).Value = "CO"
$ws.Cells.Item(17,5).Value = 10
$ws.Cells.Item(17,6).Value = 23
$ws.Cells.Item(17,7).Value = 1457
$ws.Cells.Item(17,8).Value = "Placed Actual"
$ws.Cells.Item(17,10).Value = "JOHNSTOWN"
$ws.Cells.Item(17,11).Value = "CO"
$ws.Cells.Item(17,12).Value = 170000
$ws.Cells.Item(17,13).Value = 0
$ws.Cells.Item(17,14).Value = 170000
$ws.Cells.Item(17,15).Value = "UTLX669036"

$ws.Cells.Item(18,1).Value = "TILX"
$ws.Cells.Item(18,2).Value = 253454
$ws.Cells.Item(18,3).Value = "KELIM"
$ws.Cells.Item(18,4).Value = "CO"
$ws.Cells.Item(18,5).Value = 10
$ws.Cells.Item(18,6).Value = 23
$ws.Cells.Item(18,7).Value = 2244
$ws.Cells.Item(18,8).Value = "Arrive In-Transit"
$ws.Cells.Item(18,10).Value = "JOHNSTOWN"
$ws.Cells.Item(18,11).Value = "CO"
$ws.Cells.Item(18,12).Value = 180000
$ws.Cells.Item(18,13).Value = 0
$ws.Cells.Item(18,14).Value = 180000
$ws.Cells.Item(18,15).Value = "TILX253454"

$ws.Cells.Item(19,1).Value = "UTLX"
$ws.Cells.Item(19,2).Value = 669029
$ws.Cells.Item(19,3).Value = "KELIM"
$ws.Cells.Item(19,4).Value = "CO"
$ws.Cells.Item(19,5).Value = 10
$ws.Cells.Item(19,6).Value = 23
$ws.Cells.Item(19,7).Value = 2244
$ws.Cells.Item(19,8).Value = "Arrive In-Transit"
$ws.Cells.Item(19,10).Value = "JOHNSTOWN"
$ws.Cells.Item(19,11).Value = "CO"
$ws.Cells.Item(19,12).Value = 179850
$ws.Cells.Item(19,13).Value = 0
$ws.Cells.Item(19,14).Value = 179850
$ws.Cells.Item(19,15).Value = "UTLX669029"

# Match the workbook's persisted selection (O3:O19) to the new data extent.
$ws.Range("O3:O19").Select()
